# Fix the exhibit attachment cover page placeholders:
#   users[0].address.address  ->  users[0].address.one_line
#   users[0].phone_number     ->  users[0].mobile_number

$d = $word.ActiveDocument

$d.Content.Find.Execute("address.address", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "address.one_line", 2)

$d.Content.Find.Execute("phone_number", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "mobile_number", 2)
